# Inserts two new weekly price records for "Macroferia Regional de Talca - Papa"
# at the top of the historical block (rows 253-274), pushing the existing
# 22 rows down by two positions (new total: 276 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 253 (existing data shifts down to 255..276)
$ws.Rows.Item(253).Insert()
$ws.Rows.Item(253).Insert()

# --- New row 253: Asterix, 1a (guarda), Región del Maule ---
$ws.Cells.Item(253, 1).Value = 5
$ws.Cells.Item(253, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(253, 3).Value = "Maule"
$ws.Cells.Item(253, 4).Value = 44461
$ws.Cells.Item(253, 5).Value = 7
$ws.Cells.Item(253, 6).Value = 100114001
$ws.Cells.Item(253, 7).Value = "Papa"
$ws.Cells.Item(253, 8).Value = "Asterix"
$ws.Cells.Item(253, 9).Value = "1a (guarda)"
$ws.Cells.Item(253, 10).Value = 1200
$ws.Cells.Item(253, 11).Value = 7000
$ws.Cells.Item(253, 12).Value = 7000
$ws.Cells.Item(253, 13).Value = 7000
$ws.Cells.Item(253, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(253, 15).Value = "Región del Maule"
$ws.Cells.Item(253, 16).Value = 280
$ws.Cells.Item(253, 17).Value = 25
$ws.Cells.Item(253, 18).Value = "Hortaliza"

# --- New row 254: Rodeo, 1a (guarda lavada), Región de Los Lagos ---
$ws.Cells.Item(254, 1).Value = 5
$ws.Cells.Item(254, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(254, 3).Value = "Maule"
$ws.Cells.Item(254, 4).Value = 44461
$ws.Cells.Item(254, 5).Value = 7
$ws.Cells.Item(254, 6).Value = 100114001
$ws.Cells.Item(254, 7).Value = "Papa"
$ws.Cells.Item(254, 8).Value = "Rodeo"
$ws.Cells.Item(254, 9).Value = "1a (guarda lavada)"
$ws.Cells.Item(254, 10).Value = 1200
$ws.Cells.Item(254, 11).Value = 9000
$ws.Cells.Item(254, 12).Value = 9000
$ws.Cells.Item(254, 13).Value = 9000
$ws.Cells.Item(254, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(254, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(254, 16).Value = 360
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"
